# Updates cryptos list figures (price + 1h volume) per the Jan 13 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.775.42"
$ws.Range("E2").Value = "  -6.91%  "
$ws.Range("D3").Value = "2.545.32"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D5").Value = "296.70"
$ws.Range("E5").Value = "  -4.48%  "
$ws.Range("D6").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D6").Value = "91.35"
$ws.Range("E6").Value = "  -7.01%  "
$ws.Range("E7").Value = "  -4.29%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -5.71%  "
$ws.Range("D10").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D10").Value = "35.73"
$ws.Range("E10").Value = "  -7.71%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").Value = "  -5.72%  "
$ws.Range("D13").Value = "2.936.58"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "2.543.85"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D16").Value = "0.861"
$ws.Range("E16").Value = "  -5.60%  "
$ws.Range("D17").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D17").Value = "14.08"
$ws.Range("E17").Value = "  -4.77%  "
$ws.Range("D18").Value = "42.813.82"
$ws.Range("E18").Value = "  -7.14%  "
$ws.Range("D19").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D19").Value = "6.64"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D22").Value = "72.21"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D23").Value = "259.95"
$ws.Range("E23").Value = "  -11.02%  "
$ws.Range("D24").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -4.73%  "
$ws.Range("D25").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D25").Value = "29.49"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D26").Value = "2.10"
$ws.Range("E26").Value = "  -6.11%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  -6.84%  "
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("D30").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D30").Value = "35.96"
$ws.Range("E30").Value = "  -5.81%  "
$ws.Range("D31").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D32").Value = "150.78"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D34").Value = "3.37"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D36").Value = "0.0790"
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("E37").Value = "  -6.53%  "
$ws.Range("D38").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D38").Value = "24.32"
$ws.Range("E38").Value = "  +14.92%  "
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D40").Value = "16.13"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.078.04"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D44").Value = "3.80"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("D45").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D46").Value = "84.72"
$ws.Range("E46").Value = "  -13.07%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "2.792.93"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D50").Value = "103.77"
$ws.Range("D51").NumberFormat = "@"   # keep literal text (avoid numeric auto-coercion)
$ws.Range("D51").Value = "8.65"
$ws.Range("E51").Value = "  -9.95%  "
